$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'Rp1billion'
$ws.Range("E2").Value = '  -6.39%  '

$ws.Range("D3").Value = 'Rp56.02million'
$ws.Range("E3").Value = '  -4.37%  '

$ws.Range("D4").Value = 'Rp15.767.24'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = 'Rp6.24million'
$ws.Range("E5").Value = '  -6.07%  '

$ws.Range("D6").Value = 'Rp1.93million'
$ws.Range("E6").Value = '  -5.55%  '

$ws.Range("D7").Value = 'Rp55.87million'
$ws.Range("E7").Value = '  -4.35%  '

$ws.Range("D8").Value = 'Rp9.243.88'
$ws.Range("E8").Value = '  -9.12%  '

$ws.Range("D9").Value = 'Rp15.737.14'
$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("D10").Value = 'Rp10.714.13'
$ws.Range("E10").Value = '  -11.20%  '

$ws.Range("D11").Value = 'Rp2.385.09'
$ws.Range("E11").Value = '  -18.60%  '

$ws.Range("D12").Value = 'Rp0.516'
$ws.Range("E12").Value = '  -17.48%  '

$ws.Range("D13").Value = 'Rp611.948.73'
$ws.Range("E13").Value = '  -8.34%  '

$ws.Range("D14").Value = 'Rp64.80million'
$ws.Range("E14").Value = '  -4.15%  '

$ws.Range("D15").Value = 'Rp144.384.91'
$ws.Range("E15").Value = '  -6.98%  '

$ws.Range("D16").Value = 'Rp2.142.22'
$ws.Range("E16").Value = '  -2.92%  '

$ws.Range("D17").Value = 'Rp216.818.00'
$ws.Range("E17").Value = '  +7.89%  '

$ws.Range("D18").Value = 'Rp55.98million'
$ws.Range("E18").Value = '  -4.07%  '

$ws.Range("D19").Value = 'Rp294.305.30'
$ws.Range("E19").Value = '  -8.70%  '

$ws.Range("D20").Value = 'Rp1billion'
$ws.Range("E20").Value = '  -6.30%  '

$ws.Range("D21").Value = 'Rp16.003.21'
$ws.Range("E21").Value = '  -10.10%  '

$ws.Range("D22").Value = 'Rp6.20million'
$ws.Range("E22").Value = '  -12.62%  '

$ws.Range("D23").Value = 'Rp219.093.97'
$ws.Range("E23").Value = '  -2.57%  '

$ws.Range("D24").Value = 'Rp1.29million'
$ws.Range("E24").Value = '  -8.42%  '

$ws.Range("D25").Value = 'Rp45.854.27'
$ws.Range("E25").Value = '  -5.77%  '

$ws.Range("D26").Value = 'Rp85.496.39'
$ws.Range("E26").Value = '  +9.04%  '

$ws.Range("D27").Value = 'Rp533.626.91'
$ws.Range("E27").Value = '  -11.76%  '

$ws.Range("D28").Value = 'Rp47.314.96'
$ws.Range("E28").Value = '  -7.89%  '

$ws.Range("D29").Value = 'Rp137.593.61'
$ws.Range("E29").Value = '  -15.22%  '

$ws.Range("D30").Value = 'Rp188.102.82'
$ws.Range("E30").Value = '  -3.06%  '

$ws.Range("D31").Value = 'Rp40.855.38'
$ws.Range("E31").Value = '  -7.34%  '

$ws.Range("D32").Value = 'Rp1.767.26'
$ws.Range("E32").Value = '  -5.62%  '

$ws.Range("D33").Value = 'Rp106.456.81'
$ws.Range("E33").Value = '  -5.19%  '

$ws.Range("D34").Value = 'Rp2.336.72'
$ws.Range("E34").Value = '  -7.07%  '

$ws.Range("D35").Value = 'Rp15.732.56'

$ws.Range("D36").Value = 'Rp577.016.21'
$ws.Range("E36").Value = '  -8.31%  '

$ws.Range("D37").Value = 'Rp844.556.17'
$ws.Range("E37").Value = '  -4.61%  '

$ws.Range("D38").Value = 'Rp688.67'
$ws.Range("E38").Value = '  -10.19%  '

$ws.Range("D39").Value = 'Rp15.689.07'
$ws.Range("E39").Value = '  -0.20%  '

$ws.Range("D40").Value = 'Rp0.104'
$ws.Range("E40").Value = '  -11.72%  '

$ws.Range("D41").Value = 'Rp42.108.45'
$ws.Range("E41").Value = '  -11.66%  '

$ws.Range("D42").Value = 'Rp2.053.16'
$ws.Range("E42").Value = '  -11.00%  '

$ws.Range("D43").Value = 'Rp48.394.95'
$ws.Range("E43").Value = '  +16.05%  '

$ws.Range("D44").Value = 'Rp2.23million'
$ws.Range("E44").Value = '  -4.19%  '

$ws.Range("D45").Value = 'Rp408.351.54'
$ws.Range("E45").Value = '  +1.97%  '

$ws.Range("D46").Value = 'Rp30.825.08'
$ws.Range("E46").Value = '  -5.63%  '

$ws.Range("D47").Value = 'Rp48.448.86'
$ws.Range("E47").Value = '  -9.99%  '

$ws.Range("D48").Value = 'Rp63.987.90'
$ws.Range("E48").Value = '  -5.00%  '

$ws.Range("D49").Value = 'Rp38.798.67'
$ws.Range("E49").Value = '  -8.24%  '

$ws.Range("D50").Value = 'Rp41.607.07'
$ws.Range("E50").Value = '  -8.93%  '

$ws.Range("D51").Value = 'Rp4.346.80'
$ws.Range("E51").Value = '  -9.05%  '
